# "Avances en ABM de usuarios"
# Adds two new worksheets ("Niveles de Usuarios" and "Usuarios") with their
# data + Excel Tables, and updates the tab/selection state on the existing
# sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) New sheet: "Niveles de Usuarios" (created first so it gets sheetId 3,
#    matching the history implied by the target workbook: Usuarios ends up
#    with sheetId 4 despite being moved before this sheet in tab order).
# ---------------------------------------------------------------------
$wsNiveles = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsNiveles.Name = "Niveles de Usuarios"

$wsNiveles.Range("A1").Value = "Numero Nivel"
$wsNiveles.Range("B1").Value = "Nombre"
$wsNiveles.Range("A2").Value = 0
$wsNiveles.Range("B2").Value = "Admin"
$wsNiveles.Range("A3").Value = 1
$wsNiveles.Range("B3").Value = "Gerente"
$wsNiveles.Range("A4").Value = 2
$wsNiveles.Range("B4").Value = "Mesero"
$wsNiveles.Range("A5").Value = "Escalable, agregando más niveles"

$rngNiveles = $wsNiveles.Range("A1:B5")
$rngNiveles.HorizontalAlignment = -4108
$rngNiveles.VerticalAlignment = -4108

$loNiveles = $wsNiveles.ListObjects.Add(1, $wsNiveles.Range("A1:B5"), $null, 1)
$loNiveles.Name = "Tabla1"
$loNiveles.TableStyle = "TableStyleMedium10"

$wsNiveles.Columns.Item(1).ColumnWidth = 29.877604166666668

# ---------------------------------------------------------------------
# 2) New sheet: "Usuarios" (created second, so it gets sheetId 4; it is
#    then moved before "Niveles de Usuarios" to match the target tab order).
# ---------------------------------------------------------------------
$wsUsuarios = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsUsuarios.Name = "Usuarios"

$wsUsuarios.Range("A1").Value = "Usuario"
$wsUsuarios.Range("B1").Value = "Contraseña"
$wsUsuarios.Range("C1").Value = "Tipo Perfil"
$wsUsuarios.Range("D1").Value = "idNivel"

$wsUsuarios.Range("A2").Value = "Admin"
$wsUsuarios.Range("B2").Value = "Admin123!"
$wsUsuarios.Range("C2").Value = "Administrador"
$wsUsuarios.Range("D2").Value = 0

$wsUsuarios.Range("A3").Value = 123456
$wsUsuarios.Range("B3").Value = "Gerente123!"
$wsUsuarios.Range("C3").Value = "Gerente"
$wsUsuarios.Range("D3").Value = 1

$wsUsuarios.Range("A4").Value = 654321
$wsUsuarios.Range("B4").Value = "Mesero123!"
$wsUsuarios.Range("C4").Value = "Mesero"
$wsUsuarios.Range("D4").Value = 2

$rngUsuarios = $wsUsuarios.Range("A1:D4")
$rngUsuarios.HorizontalAlignment = -4108
$rngUsuarios.VerticalAlignment = -4108

$loUsuarios = $wsUsuarios.ListObjects.Add(1, $wsUsuarios.Range("A1:D4"), $null, 1)
$loUsuarios.Name = "Tabla3"
$loUsuarios.TableStyle = "TableStyleMedium10"

$wsUsuarios.Columns.Item(1).ColumnWidth = 17.307291666666668
$wsUsuarios.Columns.Item(2).ColumnWidth = 12.307291666666666
$wsUsuarios.Columns.Item(3).ColumnWidth = 12.877604166666666

# Move "Usuarios" so it sits before "Niveles de Usuarios" in the tab order.
$wsUsuarios.Move($wsNiveles)

# Re-fetch the worksheet handles by name after Move() since the old
# object references can no longer be trusted to point at the same sheet.
$wsNiveles = $wb.Worksheets.Item("Niveles de Usuarios")
$wsUsuarios = $wb.Worksheets.Item("Usuarios")

# ---------------------------------------------------------------------
# 3) Selections / active sheet, matching the target sheetViews.
# ---------------------------------------------------------------------
$wsNiveles.Range("C3").Select()
$wsUsuarios.Range("D10").Select()

$wsTareas = $wb.Worksheets.Item("Tareas divididas")
$wsTareas.Range("B1").Select()

$wsHitos = $wb.Worksheets.Item("Hitos")
$wsHitos.Activate()
$wsHitos.Range("B5").Select()
